$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-69 (Generation 0-67): Fitness column C changes from 7573 to 7310
$ws.Range("C2:C69").Value = 7310

# Rows 70-252 (Generation 68-250): Fitness column C changes from 7573 to 7293
$ws.Range("C70:C252").Value = 7293
